# Update stats for 2025-09 (row 22 of Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B22").Value = 6292
$ws.Range("D22").Value = 5818485
$ws.Range("E22").Value = 924.7433248569612
$ws.Range("F22").Value = 8.314684110862448
$ws.Range("H22").Value = 26.53391643977418
